# Word COM-interop script
#
# Applies the diff: removes spell/grammar proof-error markup (which forced
# runs to be split around "Solr"/"Lucene"/"RavenDb"/etc.) by merging the
# surrounding runs back into single <w:r> runs, adds a new "Paged resources"
# bullet (moving the _GoBack bookmark onto it), and removes the old stray
# bookmark paragraph that used to sit just above the "Phase 1" Heading3.
#
# Strategy: Word's InsertXML on a Range replaces exactly that range's
# content with the supplied OOXML, so each paragraph we need to simplify is
# rewritten wholesale with clean markup (no w:proofErr) using the paragraph's
# own Range. For paragraph insertion/removal we expand the target Range to
# cover the neighbouring paragraph(s) so the paragraph-mark count changes by
# the right amount, and resupply the content of any paragraph that must
# survive unchanged.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml {
    param(
        [int]$Index,
        [string]$InnerXml
    )
    $p = $d.Paragraphs.Item($Index)
    $xml = "<w:p $wNs>$InnerXml</w:p>"
    $p.Range.InsertXML($xml)
}

# 1) Heading1: "Open Rasta / Solr Lucene"
Set-ParaXml 1 '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Open Rasta / Solr Lucene</w:t></w:r>'

# 2) "Create a project using Open Rasta & Solr Lucene as a prelude ..."
Set-ParaXml 3 '<w:r><w:t xml:space="preserve">Create a project using Open Rasta &amp; Solr Lucene as a prelude for my adventures at 7Digital. If there is a chance to incorporate messaging and service buses then go ahead with that too.</w:t></w:r>'

# 3) "Json"
Set-ParaXml 8 '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Json</w:t></w:r>'

# 4) "Use RavenDb and query it with lucene"
Set-ParaXml 12 '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Use RavenDb and query it with lucene</w:t></w:r>'

# 5) "Have a reporting section that can layer on some really complex lucene features"
Set-ParaXml 14 '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Have a reporting section that can layer on some really complex lucene features</w:t></w:r>'

# 6) "Etags"
Set-ParaXml 19 '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Etags</w:t></w:r>'

# 7) "Use header values and url to determin content type"
Set-ParaXml 20 '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Use header values and url to determin content type</w:t></w:r>'

# 8) "Do a StackOverflow clone - use their API"
Set-ParaXml 25 '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Do a StackOverflow clone &#8211; use their API</w:t></w:r>'

# 9) "Multiple multiple content types" (keeps the trailing " - read scenarios" run)
Set-ParaXml 38 '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Multiple multiple content types</w:t></w:r><w:r><w:t xml:space="preserve"> &#8211; read scenarios</w:t></w:r>'

# 10) "Use features that enable Hyperlinking ... resources" + "d" (two runs remain)
Set-ParaXml 42 '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Use features that enable Hyperlinking &#8211; particularly between resources</w:t></w:r><w:r><w:t>d</w:t></w:r>'

# 11+12) "WADL ?? (Check book to see what they do)" gets cleaned, and a brand new
# "Paged resources" bullet (carrying the _GoBack bookmark) is inserted right after it.
$p44 = $d.Paragraphs.Item(44)
$rng = $d.Range($p44.Range.Start, $p44.Range.End)
$xml = "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""1""/><w:numId w:val=""3""/></w:numPr></w:pPr><w:r><w:t>WADL ?? (Check book to see what they do)</w:t></w:r></w:p>" + `
       "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""1""/><w:numId w:val=""3""/></w:numPr></w:pPr><w:r><w:t>Paged resources</w:t></w:r><w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/></w:p>"
$rng.InsertXML($xml)

# 13) "Phase 3 - Enter Lucene"
Set-ParaXml 47 '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Phase 3 &#8211; </w:t></w:r><w:r><w:t>Enter Lucene</w:t></w:r>'

# 14) "Allow filtering on entities - facilitated in the back end by Lucen"
Set-ParaXml 48 '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Allow filtering on entities &#8211; facilitated in the back end by Lucen</w:t></w:r>'

# 15) "Maybe add a reporting feature to expand on the usage of Solr Lucene?"
Set-ParaXml 50 '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Maybe add a reporting feature to expand on the usage of Solr Lucene?</w:t></w:r>'

# 16) Remove the stray bookmark paragraph + trailing empty paragraph that used to sit
# right before the "Phase 1" Heading3 (re-supplying the preceding "Look at more doc"
# paragraph's own content so it survives unchanged).
$p56 = $d.Paragraphs.Item(57)
$p58 = $d.Paragraphs.Item(59)
$rng2 = $d.Range($p56.Range.Start, $p58.Range.End)
$xml2 = "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""1""/><w:numId w:val=""3""/></w:numPr></w:pPr><w:r><w:t>Look at more doc</w:t></w:r></w:p>"
$rng2.InsertXML($xml2)

# 17) "Grab some data from the api and push that into RavenDb - check it all still works"
# (index shifted by +1 after the "Paged resources" insertion and -2 after the bookmark
# paragraph removal, net -1 relative to the original document)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Grab some data from the*") {
        Set-ParaXml $i '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Grab some data from the api and push that into RavenDb &#8211; check it all still works</w:t></w:r>'
        break
    }
}
